$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "B2" = 1.02
    "C2" = 1.021323551346859
    "D2" = 1.023091967653067
    "E2" = 1.022210100574782
    "F2" = 1.019752484637941
    "I2" = 1.028395850748849
    "J2" = 1.026515664652551
    "K2" = 1.025924423788299
    "L2" = 1.02504515356685
    "M2" = 1.022594800432511
    "N2" = 1.027973434639351
    "B3" = 1.02
    "C3" = 1.022333498300389
    "D3" = 1.023964303143492
    "E3" = 1.023068577542446
    "F3" = 1.02140221847953
    "I3" = 1.028507788645528
    "J3" = 1.027162866277341
    "K3" = 1.026603385778147
    "L3" = 1.025710106318863
    "M3" = 1.024048310287277
    "N3" = 1.028621555364682
    "B4" = 1.02
    "C4" = 1.02298676628891
    "D4" = 1.024528846358328
    "E4" = 1.023624246355038
    "F4" = 1.022469335143573
    "I4" = 1.028578569745956
    "J4" = 1.027580858730693
    "K4" = 1.027042170811606
    "L4" = 1.026139916284943
    "M4" = 1.024988006075937
    "N4" = 1.029040141415263
    "B5" = 1.02
    "C5" = 1.023261344668148
    "D5" = 1.024766200576589
    "E5" = 1.02385789193929
    "F5" = 1.022917868865606
    "I5" = 1.028607930943775
    "J5" = 1.02775639401233
    "K5" = 1.027226505000477
    "L5" = 1.026320498793021
    "M5" = 1.025382863181545
    "N5" = 1.029215925977134
    "B6" = 1.02
    "C6" = 1.023307444356682
    "D6" = 1.024806054586947
    "E6" = 1.023897124550611
    "F6" = 1.022993175032694
    "I6" = 1.028612837641335
    "J6" = 1.027785856103245
    "K6" = 1.0272574478474
    "L6" = 1.026350812971254
    "M6" = 1.025449150396186
    "N6" = 1.029245429907597
    "B7" = 1.02
    "C7" = 1.022990435435641
    "D7" = 1.024532017817231
    "E7" = 1.023627368172021
    "F7" = 1.022475328787378
    "I7" = 1.028578963624643
    "J7" = 1.027583204982823
    "K7" = 1.027044634408554
    "L7" = 1.026142329667603
    "M7" = 1.024993282919177
    "N7" = 1.02904249099934
    "B8" = 1.02
    "C8" = 1.021664916854687
    "D8" = 1.023386760274167
    "E8" = 1.02250019033091
    "F8" = 1.02031010092345
    "I8" = 1.028434022302896
    "J8" = 1.026734553375457
    "K8" = 1.026153996005481
    "L8" = 1.025269972649958
    "M8" = 1.02308619467912
    "N8" = 1.028192634209349
    "B9" = 1.02
    "C9" = 1.019327344319039
    "D9" = 1.021369302963189
    "E9" = 1.020515308285222
    "F9" = 1.016491543113879
    "I9" = 1.028165983295051
    "J9" = 1.025233047466533
    "K9" = 1.024580358450087
    "L9" = 1.023729240569205
    "M9" = 1.019719101481594
    "N9" = 1.02668899598981
    "B10" = 1.02
    "C10" = 1.017767669775138
    "D10" = 1.020024740763928
    "E10" = 1.019192950747341
    "F10" = 1.013943292954748
    "I10" = 1.027978801474274
    "J10" = 1.024227928478905
    "K10" = 1.023528405197058
    "L10" = 1.022699693652642
    "M10" = 1.017469593951103
    "N10" = 1.02568244961793
    "B11" = 1.02
    "C11" = 1.017091993371728
    "D11" = 1.019442622967475
    "E11" = 1.018620564871189
    "F11" = 1.012839173471416
    "I11" = 1.027895736853946
    "J11" = 1.023791716962489
    "K11" = 1.023072212692453
    "L11" = 1.022253315683496
    "M11" = 1.016494316597909
    "N11" = 1.025245618631127
    "B12" = 1.02
    "C12" = 1.016840966666676
    "D12" = 1.019226410991181
    "E12" = 1.018407985482938
    "F12" = 1.012428939336233
    "I12" = 1.027864580419826
    "J12" = 1.023629539319205
    "K12" = 1.022902658200305
    "L12" = 1.022087423730511
    "M12" = 1.016131863851906
    "N12" = 1.025083210676987
    "B13" = 1.02
    "C13" = 1.016894815048468
    "D13" = 1.019272788636985
    "E13" = 1.018453583122527
    "F13" = 1.012516941288066
    "I13" = 1.02787127727123
    "J13" = 1.023664333741247
    "K13" = 1.022939032944035
    "L13" = 1.022123012077156
    "M13" = 1.016209619968941
    "N13" = 1.0251180545111
    "B14" = 1.02
    "C14" = 1.017071244461603
    "D14" = 1.019424750566194
    "E14" = 1.018602992378611
    "F14" = 1.012805265783787
    "I14" = 1.02789316762507
    "J14" = 1.023778314357351
    "K14" = 1.023058199396391
    "L14" = 1.022239604787077
    "M14" = 1.016464360085122
    "N14" = 1.025232196992752
    "B15" = 1.02
    "C15" = 1.017179941728671
    "D15" = 1.019518380966175
    "E15" = 1.018695052343982
    "F15" = 1.012982896534796
    "I15" = 1.027906614905369
    "J15" = 1.023848521774547
    "K15" = 1.023131607942452
    "L15" = 1.022311429818691
    "M15" = 1.016621288312594
    "N15" = 1.025302504112534
    "B16" = 1.02
    "C16" = 1.01781250519426
    "D16" = 1.020063375784741
    "E16" = 1.019230942371072
    "F16" = 1.014016554056878
    "I16" = 1.027984271773658
    "J16" = 1.024256857511186
    "K16" = 1.023558666602647
    "L16" = 1.022729306076598
    "M16" = 1.017534293461831
    "N16" = 1.025711419732756
    "B17" = 1.02
    "C17" = 1.018209207094361
    "D17" = 1.020405259444424
    "E17" = 1.019567146284789
    "F17" = 1.014664744104371
    "I17" = 1.028032444759113
    "J17" = 1.024512730575755
    "K17" = 1.023826364083041
    "L17" = 1.022991273846658
    "M17" = 1.018106664263741
    "N17" = 1.025967656166425
    "B18" = 1.02
    "C18" = 1.018440565132888
    "D18" = 1.020604682780258
    "E18" = 1.019763268016744
    "F18" = 1.01504275441167
    "I18" = 1.028060348958484
    "J18" = 1.024661881695149
    "K18" = 1.023982440982259
    "L18" = 1.023144019545911
    "M18" = 1.018440400683941
    "N18" = 1.026117019097517
    "B19" = 1.02
    "C19" = 1.018519446946162
    "D19" = 1.020672682421014
    "E19" = 1.019830143808276
    "F19" = 1.015171634878181
    "I19" = 1.028069830614881
    "J19" = 1.024712722229573
    "K19" = 1.024035647909487
    "L19" = 1.023196092456479
    "M19" = 1.018554176417558
    "N19" = 1.026167931831331
    "B20" = 1.02
    "C20" = 1.018166647998657
    "D20" = 1.020368577699824
    "E20" = 1.019531072753954
    "F20" = 1.014595206537137
    "I20" = 1.028027296349786
    "J20" = 1.024485287681255
    "K20" = 1.023797649561009
    "L20" = 1.022963172946121
    "M20" = 1.01804526652749
    "N20" = 1.025940174299867
    "B21" = 1.02
    "C21" = 1.017019291798491
    "D21" = 1.019380001224864
    "E21" = 1.018558994245795
    "F21" = 1.012720364647295
    "I21" = 1.027886729813003
    "J21" = 1.023744754033476
    "K21" = 1.023023110730327
    "L21" = 1.022205273554452
    "M21" = 1.016389350808983
    "N21" = 1.025198589009368
    "B22" = 1.02
    "C22" = 1.016297611155948
    "D22" = 1.018758515975297
    "E22" = 1.017947984501273
    "F22" = 1.011540907664628
    "I22" = 1.027796599668072
    "J22" = 1.023278287152531
    "K22" = 1.022535523742308
    "L22" = 1.021728246749576
    "M22" = 1.015347100424416
    "N22" = 1.024731459691948
    "B23" = 1.02
    "C23" = 1.016680215849664
    "D23" = 1.019087970410516
    "E23" = 1.018271875910224
    "F23" = 1.012166226361901
    "I23" = 1.027844545282656
    "J23" = 1.023525652229012
    "K23" = 1.022794060271915
    "L23" = 1.02198117577764
    "M23" = 1.015899724732244
    "N23" = 1.02497917605521
    "B24" = 1.02
    "C24" = 1.018185878718888
    "D24" = 1.020385152575112
    "E24" = 1.019547372768746
    "F24" = 1.014626627778686
    "I24" = 1.028029623294602
    "J24" = 1.024497688238215
    "K24" = 1.023810624622758
    "L24" = 1.022975870705248
    "M24" = 1.01807300988391
    "N24" = 1.025952592467041
    "B25" = 1.02
    "C25" = 1.019931886576375
    "D25" = 1.021890790662977
    "E25" = 1.02102828896257
    "F25" = 1.017479148963601
    "I25" = 1.028236774437087
    "J25" = 1.025621945682079
    "K25" = 1.024987683975966
    "L25" = 1.024127976847932
    "M25" = 1.020590390978256
    "N25" = 1.027078446485429
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}
